# Auto-generated: updates cryptos list price/volume(1h) data per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.987.14'
$ws.Range("E2").Value = '  -3.35%  '
$ws.Range("D3").Value = '3.333.86'
$ws.Range("E3").Value = '  -5.17%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '''552.96'
$ws.Range("E5").Value = '  -4.26%  '
$ws.Range("D6").Value = '''172.50'
$ws.Range("E6").Value = '  -3.89%  '
$ws.Range("D7").Value = '''0.613'
$ws.Range("E7").Value = '  -3.86%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '3.321.33'
$ws.Range("E9").Value = '  -5.35%  '
$ws.Range("D10").Value = '''0.621'
$ws.Range("E10").Value = '  -2.36%  '
$ws.Range("E11").Value = '  +2.50%  '
$ws.Range("D12").Value = '''53.05'
$ws.Range("E12").Value = '  -5.07%  '
$ws.Range("E13").Value = '  -1.25%  '
$ws.Range("E14").Value = '  -3.22%  '
$ws.Range("D15").Value = '3.862.19'
$ws.Range("E15").Value = '  -5.28%  '
$ws.Range("D16").Value = '''18.24'
$ws.Range("E16").Value = '  -1.18%  '
$ws.Range("E17").Value = '  -3.71%  '
$ws.Range("D18").Value = '3.332.40'
$ws.Range("E18").Value = '  -5.09%  '
$ws.Range("D19").Value = '''11.81'
$ws.Range("E19").Value = '  -2.39%  '
$ws.Range("D20").Value = '63.852.16'
$ws.Range("E20").Value = '  -3.50%  '
$ws.Range("D21").Value = '''0.971'
$ws.Range("E21").Value = '  -3.80%  '
$ws.Range("D22").Value = '''427.15'
$ws.Range("E22").Value = '  +2.76%  '
$ws.Range("E23").Value = '  +8.16%  '
$ws.Range("E24").Value = '  -3.44%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = '''84.11'
$ws.Range("E25").Value = '  -2.19%  '
$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").Value = '''13.42'
$ws.Range("E26").Value = '  +0.20%  '
$ws.Range("E27").Value = '  -3.81%  '
$ws.Range("E28").Value = '  -1.90%  '
$ws.Range("E29").Value = '  -6.81%  '
$ws.Range("D30").Value = '''29.65'
$ws.Range("E30").Value = '  -3.01%  '
$ws.Range("D31").Value = '''6.67'
$ws.Range("E31").Value = '  +0.60%  '
$ws.Range("D32").Value = '''593.44'
$ws.Range("E32").Value = '  -6.19%  '
$ws.Range("D33").Value = '''11.38'
$ws.Range("E33").Value = '  -2.95%  '
$ws.Range("E34").Value = '  -3.69%  '
$ws.Range("D35").Value = '''58.13'
$ws.Range("D36").Value = '''0.999'
$ws.Range("E36").Value = '  -0.16%  '
$ws.Range("E37").Value = '  -8.33%  '
$ws.Range("D38").Value = '''3.48'
$ws.Range("E38").Value = '  +1.03%  '
$ws.Range("D39").Value = '''35.22'
$ws.Range("E39").Value = '  -5.59%  '
$ws.Range("D40").Value = '0.0₃0747'
$ws.Range("E40").Value = '  -7.08%  '
$ws.Range("D41").Value = '''0.363'
$ws.Range("E41").Value = '  -5.06%  '
$ws.Range("D42").Value = '3.093.83'
$ws.Range("E42").Value = '  -5.63%  '
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("E44").Value = '  -5.09%  '
$ws.Range("D45").Value = '''0.0405'
$ws.Range("E45").Value = '  -3.53%  '
$ws.Range("D46").Value = '''3.18'
$ws.Range("E46").Value = '  -4.01%  '
$ws.Range("E47").Value = '  -4.19%  '
$ws.Range("E48").Value = '  -2.77%  '
$ws.Range("D49").Value = '''2.58'
$ws.Range("E49").Value = '  -4.42%  '
$ws.Range("D50").Value = '''8.14'
$ws.Range("E50").Value = '  -6.19%  '
$ws.Range("D51").Value = '''132.23'
$ws.Range("E51").Value = '  -4.86%  '
